$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '27.377.63'
$ws.Cells.Item(2,5).Value = '  -0.80%  '

# Row 3
$ws.Cells.Item(3,4).Value = '1.637.36'
$ws.Cells.Item(3,5).Value = '  -1.70%  '

# Row 4
$ws.Cells.Item(4,5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '210.91'
$ws.Cells.Item(5,5).Value = '  -1.85%  '

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.529'
$ws.Cells.Item(6,5).Value = '  +3.55%  '

# Row 7
$ws.Cells.Item(7,5).Value = '  +0.05%  '

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '23.07'
$ws.Cells.Item(8,5).Value = '  -2.40%  '

# Row 9
$ws.Cells.Item(9,5).Value = '  -2.88%  '

# Row 10
$ws.Cells.Item(10,5).Value = '  -2.22%  '

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.0888'
$ws.Cells.Item(11,5).Value = '  +1.08%  '

# Row 12
$ws.Cells.Item(12,4).Value = '1.870.45'
$ws.Cells.Item(12,5).Value = '  -1.62%  '

# Row 13
$ws.Cells.Item(13,4).Value = '1.638.10'
$ws.Cells.Item(13,5).Value = '  -2.73%  '

# Row 14
$ws.Cells.Item(14,5).Value = '  -3.07%  '

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.556'
$ws.Cells.Item(15,5).Value = '  -1.17%  '

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '64.20'
$ws.Cells.Item(16,5).Value = '  -3.00%  '

# Row 17
$ws.Cells.Item(17,4).Value = '27.351.80'
$ws.Cells.Item(17,5).Value = '  -0.90%  '

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '229.78'
$ws.Cells.Item(18,5).Value = '  -5.26%  '

# Row 19
$ws.Cells.Item(19,5).Value = '  -1.40%  '

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '7.54'
$ws.Cells.Item(20,5).Value = '  -0.57%  '

# Row 21
$ws.Cells.Item(21,5).Value = '  +0.04%  '

# Row 22
$ws.Cells.Item(22,5).Value = '  -3.97%  '

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '9.39'
$ws.Cells.Item(23,5).Value = '  +1.24%  '

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '2.02'
$ws.Cells.Item(24,5).Value = '  -0.85%  '

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '147.94'
$ws.Cells.Item(25,5).Value = '  +0.85%  '

# Row 26
$ws.Cells.Item(26,5).Value = '  -3.57%  '

# Row 27
$ws.Cells.Item(27,5).Value = '  +1.24%  '

# Row 28
$ws.Cells.Item(28,5).Value = '  -0.03%  '

# Row 29
$ws.Cells.Item(29,5).Value = '  -5.54%  '

# Row 30
$ws.Cells.Item(30,5).Value = '  -4.69%  '

# Row 31
$ws.Cells.Item(31,5).Value = '  -3.59%  '

# Row 32
$ws.Cells.Item(32,5).Value = '  -2.46%  '

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '3.11'
$ws.Cells.Item(33,5).Value = '  -0.10%  '

# Row 34
$ws.Cells.Item(34,4).Value = '1.405.89'
$ws.Cells.Item(34,5).Value = '  -4.33%  '

# Row 35
$ws.Cells.Item(35,5).Value = '  -0.06%  '

# Row 36
$ws.Cells.Item(36,5).Value = '  -0.19%  '

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.561'
$ws.Cells.Item(37,5).Value = '  -2.09%  '

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.877'
$ws.Cells.Item(38,5).Value = '  -5.72%  '

# Row 39
$ws.Cells.Item(39,5).Value = '  -3.63%  '

# Row 40
$ws.Cells.Item(40,5).Value = '  +0.61%  '

# Row 41
$ws.Cells.Item(41,5).Value = '  +0.03%  '

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '2.46'
$ws.Cells.Item(42,5).Value = '  -1.92%  '

# Row 43
$ws.Cells.Item(43,2).Value = 'FraxShare'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '5.46'
$ws.Cells.Item(43,5).Value = '  +1.00%  '

# Row 44
$ws.Cells.Item(44,2).Value = 'MXToken'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '2.23'
$ws.Cells.Item(44,5).Value = '  +0.58%  '

# Row 45
$ws.Cells.Item(45,2).Value = 'TrustWalletToken'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.789'
$ws.Cells.Item(45,5).Value = '  +0.36%  '

# Row 46
$ws.Cells.Item(46,2).Value = 'Aave'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '64.36'
$ws.Cells.Item(46,5).Value = '  -7.22%  '

# Row 47
$ws.Cells.Item(47,4).Value = '1.779.52'
$ws.Cells.Item(47,5).Value = '  -1.62%  '

# Row 48
$ws.Cells.Item(48,5).Value = '  -4.60%  '

# Row 49
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '87.15'
$ws.Cells.Item(49,5).Value = '  -2.52%  '

# Row 50
$ws.Cells.Item(50,5).Value = '  -2.77%  '

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.0986'
$ws.Cells.Item(51,5).Value = '  -4.01%  '
